$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that receive numeric-looking text values need to be forced to Text
# format before assignment so Excel stores them as strings (matching the source data),
# not auto-converted numbers. We do this with a batch NumberFormat on column D, write all
# D/B/C/E values, then reset the number format back to General so no stray formatting
# is left behind on the cells.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "57.240.33"
$ws.Range("E2").Value = "  +5.19%  "

$ws.Range("D3").Value = "2.359.04"
$ws.Range("E3").Value = "  +3.31%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").Value = "519.41"
$ws.Range("E5").Value = "  +3.60%  "

$ws.Range("D6").Value = "134.67"
$ws.Range("E6").Value = "  +3.69%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +2.21%  "

$ws.Range("D9").Value = "2.354.70"
$ws.Range("E9").Value = "  +2.74%  "

$ws.Range("E10").Value = "  +7.86%  "

$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +7.05%  "

$ws.Range("E13").Value = "  +2.37%  "

$ws.Range("D14").Value = "23.93"
$ws.Range("E14").Value = "  +3.71%  "

$ws.Range("D15").Value = "2.754.18"
$ws.Range("E15").Value = "  +2.28%  "

$ws.Range("D16").Value = "57.064.79"
$ws.Range("E16").Value = "  +4.72%  "

$ws.Range("E17").Value = "  +3.66%  "

$ws.Range("D18").Value = "2.344.67"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("E19").Value = "  +2.53%  "

$ws.Range("D20").Value = "4.31"
$ws.Range("E20").Value = "  +3.62%  "

$ws.Range("D21").Value = "321.27"
$ws.Range("E21").Value = "  +5.18%  "

$ws.Range("E22").Value = "  +6.45%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "61.30"
$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "0.161"
$ws.Range("E25").Value = "  +7.52%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +5.81%  "

$ws.Range("D28").Value = "171.87"
$ws.Range("E28").Value = "  -0.46%  "

$ws.Range("E29").Value = "  +10.42%  "

$ws.Range("D30").Value = "0.0₃0741"
$ws.Range("E30").Value = "  +4.69%  "

$ws.Range("E31").Value = "  +4.61%  "

$ws.Range("E32").Value = "  +3.85%  "

$ws.Range("D33").Value = "18.39"
$ws.Range("E33").Value = "  +2.54%  "

$ws.Range("D35").Value = "0.966"
$ws.Range("E35").Value = "  +3.48%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("E37").Value = "  +4.99%  "

$ws.Range("E38").Value = "  +8.04%  "

$ws.Range("D39").Value = "37.61"
$ws.Range("E39").Value = "  +4.02%  "

$ws.Range("E40").Value = "  +7.29%  "

$ws.Range("E41").Value = "  +1.91%  "

$ws.Range("D42").Value = "140.46"
$ws.Range("E42").Value = "  +12.42%  "

$ws.Range("E43").Value = "  +6.10%  "

$ws.Range("D44").Value = "278.46"
$ws.Range("E44").Value = "  +12.62%  "

$ws.Range("D45").Value = "5.18"
$ws.Range("E45").Value = "  +2.58%  "

$ws.Range("E46").Value = "  +3.68%  "

$ws.Range("E47").Value = "  +3.80%  "

$ws.Range("E48").Value = "  +3.02%  "

$ws.Range("E49").Value = "  +2.11%  "

$ws.Range("E50").Value = "  +4.57%  "

$ws.Range("D51").Value = "17.01"
$ws.Range("E51").Value = "  +3.18%  "

# Restore column D to its default style (removes the temporary Text number format,
# while keeping the values stored as text).
$colD.Style = "Normal"
